$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the rows for dcdb0ae3... and 571c7103... swap positions
# (dcdb0ae3... now listed before 571c7103...) and the status for
# 197fab40... / dcdb0ae3... flips from "Ready for handoff" to
# "In Translation".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"

$wsOverview.Range("A4").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"
$wsOverview.Range("D4").Value = "2016-03-21 22:36:37"

$wsOverview.Range("A5").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-21 22:34:46"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7d9cae1f907f2f0e19d72d98df700e24046a7a28/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d707f83dbe3a5f9640cb39e58e4f1b23177cffa4/e2e/197fab40-7705-4dd6-a3a2-ec57183f1ea7.md", [Type]::Missing, [Type]::Missing, "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/571c7103-8660-4f11-8e8b-df8803d0e27d.md", [Type]::Missing, [Type]::Missing, "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d707f83dbe3a5f9640cb39e58e4f1b23177cffa4/e2e/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md", [Type]::Missing, [Type]::Missing, "571c7103-8660-4f11-8e8b-df8803d0e27d.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row swap / status update as above.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C3").Value = "In Translation"

$wsZhCn.Range("A4").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("D4").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.zh-cn.xlf"
$wsZhCn.Range("E4").Value = "2016-03-21 22:36:31"

$wsZhCn.Range("A5").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.zh-cn.xlf"
$wsZhCn.Range("E5").Value = "2016-03-21 22:34:43"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7d9cae1f907f2f0e19d72d98df700e24046a7a28/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8da87dd7580ef39e6e092d7e540a6cab6d528832/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a734c2cb4d4d89f1fdb4589337d30735c022bd58/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/01f75f6f65a8b3d6abb033821e7cfdf48a1b05a6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d707f83dbe3a5f9640cb39e58e4f1b23177cffa4/e2e/197fab40-7705-4dd6-a3a2-ec57183f1ea7.md", [Type]::Missing, [Type]::Missing, "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6da88f2a22291b1e3470d34b3ca5af9035730286/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/571c7103-8660-4f11-8e8b-df8803d0e27d.md", [Type]::Missing, [Type]::Missing, "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86c5f1c3fcf4eeb680ca7d9a65b3ee542c89a9bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d707f83dbe3a5f9640cb39e58e4f1b23177cffa4/e2e/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md", [Type]::Missing, [Type]::Missing, "571c7103-8660-4f11-8e8b-df8803d0e27d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6da88f2a22291b1e3470d34b3ca5af9035730286/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de": same row swap / status update as above.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("E3").Value = "2016-03-21 22:36:37"

$wsDeDe.Range("A4").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("D4").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.de-de.xlf"
$wsDeDe.Range("E4").Value = "2016-03-21 22:36:37"

$wsDeDe.Range("A5").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.de-de.xlf"
$wsDeDe.Range("E5").Value = "2016-03-21 22:34:46"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7d9cae1f907f2f0e19d72d98df700e24046a7a28/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86ebc92335ca81ee9d52dae18a8f7f0bd62a2f34/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/babad688c7f84e7d63a0266026211f6fdafe15bb/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e7faf877098ed28ec94e1703d5ffd504aa2d4761/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf", [Type]::Missing, [Type]::Missing, "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d707f83dbe3a5f9640cb39e58e4f1b23177cffa4/e2e/197fab40-7705-4dd6-a3a2-ec57183f1ea7.md", [Type]::Missing, [Type]::Missing, "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/634d7cba4e29b47e0c3c4dbf92e9f0fe0975d623/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.de-de.xlf", [Type]::Missing, [Type]::Missing, "197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/571c7103-8660-4f11-8e8b-df8803d0e27d.md", [Type]::Missing, [Type]::Missing, "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34ae5c474cf2361ac996412ee2a82e4e64ab8941/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.de-de.xlf", [Type]::Missing, [Type]::Missing, "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d707f83dbe3a5f9640cb39e58e4f1b23177cffa4/e2e/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md", [Type]::Missing, [Type]::Missing, "571c7103-8660-4f11-8e8b-df8803d0e27d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/634d7cba4e29b47e0c3c4dbf92e9f0fe0975d623/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.de-de.xlf", [Type]::Missing, [Type]::Missing, "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.de-de.xlf")

$wb.Save()
